$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# ------------------------------------------------------------------
# 1) Insert a new row above row 11 (this shifts old rows 11-18 down
#    to 12-19, carrying their values AND formatting with them).
# ------------------------------------------------------------------
$ws.Rows("11:11").Insert()

# ------------------------------------------------------------------
# 2) The brand new row 11 comes back blank/default-formatted. Give it
#    the same look as row 10 (immediately above), which is the
#    formatting the new "TargetTable" record row ends up with.
# ------------------------------------------------------------------
$ws.Range("A10:J10").Copy()
$ws.Range("A11:J11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Fix up the SEQ column: the rows that got pushed down keep their
#    old literal numbers (3,4,5,6,7,8,9,10) but must now read
#    4,5,6,7,8,9,10,11 since a new #3 was inserted above them.
# ------------------------------------------------------------------
$ws.Range("A12").Value = 4
$ws.Range("A13").Value = 5
$ws.Range("A14").Value = 6
$ws.Range("A15").Value = 7
$ws.Range("A16").Value = 8
$ws.Range("A17").Value = 9
$ws.Range("A18").Value = 10
$ws.Range("A19").Value = 11

# ------------------------------------------------------------------
# 4) Fill in the new row 11 ("TargetTable") content.
# ------------------------------------------------------------------
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "TargetTable"
$ws.Range("C11").Value = "對應資料表"
$ws.Range("D11").Value = "VARCHAR2"
$ws.Range("E11").Value = 30
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = "大小寫需完全符合"
$ws.Range("H11").Value = "擺NULL"
$ws.Range("I11").Value = "擺對應資料表名稱"
$ws.Range("J11").Value = ""

# ------------------------------------------------------------------
# 5) TableName row (row 10) gets a note in the remarks column.
# ------------------------------------------------------------------
$ws.Range("G10").Value = "大小寫需完全符合"

# ------------------------------------------------------------------
# 6) Match the saved selection/cursor position.
# ------------------------------------------------------------------
$ws.Range("H12").Select()
